# Insert a new product row ("LACRITEARS EYE DROPS 15 ML") right after the
# "IVY PRONT  SYRUP" row (row 29) and before the "LIBRAX 30 SUGAR COATED TAB"
# row (row 30), shifting every following row down by one, renumbering the
# sequence column (A), and updating the running total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertRow = 30
$lastOldDataRow = 55      # last data row before insertion
$lastNewDataRow = $lastOldDataRow + 1   # 56, after insertion
$oldTotalRow = 56
$newTotalRow = $oldTotalRow + 1         # 57

# 1) Insert a blank row at row 30; everything from 30..57 shifts to 31..58.
$ws.Rows.Item($insertRow).Insert()

# 2) Copy the formatting (styles, fonts, borders, number formats) of the row
#    that used to be row 30 (now row 31, the "LIBRAX..." row) onto the new
#    blank row so it reuses the existing row style instead of creating new
#    style records.
$ws.Range("A31:N31").Copy()
$ws.Range("A30:N30").PasteSpecial(-4122)  # xlPasteFormats

# 3) Recreate the merged cell groups for the new row (B:G, H:K, L:M), same
#    as every other data row.
$ws.Range("B30:G30").Merge()
$ws.Range("H30:K30").Merge()
$ws.Range("L30:M30").Merge()

# 4) The row heights in this sheet are fixed per absolute row position (they
#    do not travel with the shifted content), so restore rows 30..55 to the
#    same heights they had before the insertion, and give the brand-new rows
#    (the new last data row and the new total row) their own heights.
$rowHeights = @{
    30 = 25.5;  31 = 24.75; 32 = 25.5;  33 = 25.5;  34 = 24.75; 35 = 25.5;
    36 = 24.75; 37 = 25.5;  38 = 25.5;  39 = 24.75; 40 = 25.5;  41 = 24.75;
    42 = 25.5;  43 = 25.5;  44 = 24.75; 45 = 25.5;  46 = 24.75; 47 = 25.5;
    48 = 25.5;  49 = 24.75; 50 = 25.5;  51 = 24.75; 52 = 25.5;  53 = 25.5;
    54 = 24.75; 55 = 25.5
}
foreach ($r in $rowHeights.Keys) {
    $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
}
$ws.Rows.Item($lastNewDataRow).RowHeight = 24.75
$ws.Rows.Item($newTotalRow).RowHeight = 26.25

# 5) Fill in the values for the new row.
$ws.Range("A30").Value = 27
$ws.Range("B30").Value = "LACRITEARS EYE DROPS 15 ML"
$ws.Range("H30").Value = "1:0"
$ws.Range("L30").Value = 49
$ws.Range("N30").Value = "1:0"

# 6) Renumber the sequence column (A) for every row that shifted down, so the
#    numbering stays consecutive (28, 29, 30, ... 53).
for ($r = 31; $r -le $lastNewDataRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 3
}

# 7) Update the running total in column K (was 2580.17, now +49 = 2629.17).
$ws.Cells.Item($newTotalRow, 11).Value = 2629.1700000000001

Write-Host "Row inserted and renumbered successfully"
